$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Sheet view: scroll position + selection ---
$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("E18").Select()

# --- Row 9: Charge Amount becomes text "$1,000" ---
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "$1,000"

# --- Row 10: Charge Amount becomes text "$2,000" ---
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "$2,000"

# --- Row 11: Building Abbreviation fixed typo, and number format switched to Text (value unchanged) ---
$ws.Range("A11").Value = "SampleAbbreviation2"
$ws.Range("E11").Value = 15
$ws.Range("E11").NumberFormat = "@"

# --- Rows 12-14: number format switched to Text (values unchanged) ---
$ws.Range("E12").Value = 1100
$ws.Range("E12").NumberFormat = "@"

$ws.Range("E13").Value = 800
$ws.Range("E13").NumberFormat = "@"

$ws.Range("E14").Value = 1200
$ws.Range("E14").NumberFormat = "@"

# --- Row 15: number format switched to General (value unchanged) ---
$ws.Range("E15").Value = 1800
$ws.Range("E15").NumberFormat = "General"

# --- Row 16: number format switched to General, value updated ---
$ws.Range("E16").Value = 1900
$ws.Range("E16").NumberFormat = "General"

# --- Row 17: value updated (format stays currency) ---
$ws.Range("E17").Value = 30

# --- Row 19: explicit zero value (was blank) ---
$ws.Range("E19").Value = 0

Write-Host "done"
